$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1004
$ws.Range("F7").Value = 2517
$ws.Range("F10").Value = 901
$ws.Range("F13").Value = 1117
$ws.Range("F17").Value = 722
$ws.Range("F18").Value = 764
$ws.Range("F19").Value = 187
$ws.Range("F20").Value = 478
$ws.Range("F22").Value = 84
$ws.Range("F24").Value = 582
$ws.Range("F29").Value = 304
$ws.Range("F30").Value = 4244
$ws.Range("F36").Value = 141
$ws.Range("F37").Value = 1593
$ws.Range("F40").Value = 83
$ws.Range("F41").Value = 138
$ws.Range("F42").Value = 68

$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 11
$ws.Range("F13").Value = 13
$ws.Range("F16").Value = 182

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 724

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 724
$ws.Range("F7").Value = 1004
$ws.Range("F8").Value = 2517
$ws.Range("F11").Value = 901
$ws.Range("F14").Value = 1117
$ws.Range("F18").Value = 722
$ws.Range("F21").Value = 764
$ws.Range("F22").Value = 187
$ws.Range("F23").Value = 478
$ws.Range("F26").Value = 84
$ws.Range("F28").Value = 582
$ws.Range("F33").Value = 4244
$ws.Range("F38").Value = 141
$ws.Range("F39").Value = 1593
$ws.Range("F41").Value = 13
$ws.Range("F44").Value = 83
$ws.Range("F45").Value = 138
